$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.088.14"
$ws.Range("D3").Value = "3.757.06"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "3.754.82"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "4.384.23"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "3.766.36"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "69.090.36"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "3.901.94"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "3.719.58"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.324"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "428.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "2.813.48"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("E51").Value = "  +9.39%  "
